$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 3.3
$ws.Range("T2").Value = 11
$ws.Range("V2").Value = 9.5
$ws.Range("AD2").Value = 101
